$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-18
$values = @{
    2  = @(5, 6)
    3  = @(9, 9)
    4  = @(8, 9)
    5  = @(8, 8)
    6  = @(6, 6)
    7  = @(7, 8)
    8  = @(7, 7)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(7, 8)
    12 = @(8, 8)
    13 = @(8, 9)
    14 = @(6, 6)
    15 = @(7, 8)
    16 = @(9, 9)
    17 = @(8, 8)
    18 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
